$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Contact_Info column keeps its text formatting so values
# like "03234567890" are not coerced into numbers (losing leading zeros).
$ws.Range("C2:C6").NumberFormat = "@"

# Row 2: Ahmed -> Hamza Traders
$ws.Range("A2").Value = 16
$ws.Range("B2").Value = "Hamza Traders"
$ws.Range("C2").Value = "03234567890"

# Row 3: Gul Khan -> Hassan General Store
$ws.Range("A3").Value = 17
$ws.Range("B3").Value = "Hassan General Store"
$ws.Range("C3").Value = "03001234567"
$ws.Range("E3").Value = "Bonapapa"

# Row 4: Karim -> Rashid Mart
$ws.Range("A4").Value = 18
$ws.Range("B4").Value = "Rashid Mart"
$ws.Range("C4").Value = "03124567890"

# Row 5: Nadim -> Al Madina Store
$ws.Range("A5").Value = 19
$ws.Range("B5").Value = "Al Madina Store"
$ws.Range("C5").Value = "03451234567"

# Row 6: new row
$ws.Range("A6").Value = 20
$ws.Range("B6").Value = "Kroon Shop"
$ws.Range("C6").Value = "11243243245"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "Candyland"
